$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case "de/del/el/la/las/los/y" connector words in state/municipality names ---
$ws.Range("B24").Value = 'Amatenango De La Frontera'
$ws.Range("B27").Value = 'Bejucal De Ocampo'
$ws.Range("B35").Value = 'Chiapa De Corzo'
$ws.Range("B39").Value = 'Comitán De Domínguez'
$ws.Range("B58").Value = 'Mazapa De Madero'
$ws.Range("B62").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B68").Value = 'Salto De Agua'
$ws.Range("B69").Value = 'San Cristóbal De Las Casas'
$ws.Range("B93").Value = 'Hidalgo Del Parral'
$ws.Range("A113").Value = 'Ciudad De México'
$ws.Range("B117").Value = 'Cuajimalpa De Morelos'
$ws.Range("B138").Value = 'San Juan Del Río'
$ws.Range("A144").Value = 'Estado De México'
$ws.Range("B144").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B147").Value = 'Almoloya De Alquisiras'
$ws.Range("B148").Value = 'Almoloya De Juárez'
$ws.Range("B155").Value = 'Atizapán De Zaragoza'
$ws.Range("B163").Value = 'Chapa De Mota'
$ws.Range("B167").Value = 'Coacalco De Berriozábal'
$ws.Range("B173").Value = 'Ecatepec De Morelos'
$ws.Range("B178").Value = 'Ixtapan De La Sal'
$ws.Range("B179").Value = 'Ixtapan Del Oro'
$ws.Range("B193").Value = 'Naucalpan De Juárez'
$ws.Range("B202").Value = 'San Antonio La Isla'
$ws.Range("B203").Value = 'San Felipe Del Progreso'
$ws.Range("B204").Value = 'San Martín De Las Pirámides'
$ws.Range("B215").Value = 'Tenango Del Valle'
$ws.Range("B226").Value = 'Tlalnepantla De Baz'
$ws.Range("B231").Value = 'Valle De Bravo'
$ws.Range("B232").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B233").Value = 'Villa De Allende'
$ws.Range("B234").Value = 'Villa Del Carbón'
$ws.Range("B246").Value = 'San Miguel De Allende'
$ws.Range("B247").Value = 'Apaseo El Alto'
$ws.Range("B248").Value = 'Apaseo El Grande'
$ws.Range("B252").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B256").Value = 'Jaral Del Progreso'
$ws.Range("B266").Value = 'San Diego De La Unión'
$ws.Range("B268").Value = 'San Francisco Del Rincón'
$ws.Range("B270").Value = 'San Luis De La Paz'
$ws.Range("B271").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B272").Value = 'Silao De La Victoria'
$ws.Range("B274").Value = 'Valle De Santiago'
$ws.Range("B279").Value = 'Acapulco De Juárez'
$ws.Range("B282").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B283").Value = 'Alcozauca De Guerrero'
$ws.Range("B287").Value = 'Atenango Del Río'
$ws.Range("B288").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B290").Value = 'Atoyac De Álvarez'
$ws.Range("B291").Value = 'Ayutla De Los Libres'
$ws.Range("B294").Value = 'Buenavista De Cuéllar'
$ws.Range("B295").Value = 'Chilapa De Álvarez'
$ws.Range("B296").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B297").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B302").Value = 'Coyuca De Benítez'
$ws.Range("B303").Value = 'Coyuca De Catalán'
$ws.Range("B307").Value = 'Cuetzala Del Progreso'
$ws.Range("B308").Value = 'Cutzamala De Pinzón'
$ws.Range("B314").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B315").Value = 'Iguala De La Independencia'
$ws.Range("B317").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B320").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B323").Value = 'Mártir De Cuilapan'
$ws.Range("B336").Value = 'Taxco De Alarcón'
$ws.Range("B338").Value = 'Técpan De Galeana'
$ws.Range("B340").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B342").Value = 'Tixtla De Guerrero'
$ws.Range("B346").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B347").Value = 'Tlapa De Comonfort'
$ws.Range("B359").Value = 'Agua Blanca De Iturbide'
$ws.Range("B366").Value = 'Atotonilco El Grande'
$ws.Range("B372").Value = 'Cuautepec De Hinojosa'
$ws.Range("B377").Value = 'Huasca De Ocampo'
$ws.Range("B381").Value = 'Huejutla De Reyes'
$ws.Range("B384").Value = 'Jacala De Ledezma'
$ws.Range("B390").Value = 'Mineral Del Chico'
$ws.Range("B391").Value = 'Mineral Del Monte'
$ws.Range("B392").Value = 'Mixquiahuala De Juárez'
$ws.Range("B393").Value = 'Molango De Escamilla'
$ws.Range("B395").Value = 'Nopala De Villagrán'
$ws.Range("B396").Value = 'Omitlán De Juárez'
$ws.Range("B397").Value = 'Pachuca De Soto'
$ws.Range("B399").Value = 'Progreso De Obregón'
$ws.Range("B405").Value = 'Santiago De Anaya'
$ws.Range("B409").Value = 'Tenango De Doria'
$ws.Range("B411").Value = 'Tepehuacán De Guerrero'
$ws.Range("B412").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B415").Value = 'Tezontepec De Aldama'
$ws.Range("B423").Value = 'Tula De Allende'
$ws.Range("B424").Value = 'Tulancingo De Bravo'
$ws.Range("B425").Value = 'Villa De Tezontepec'
$ws.Range("B429").Value = 'Zacualtipán De Ángeles'
$ws.Range("B430").Value = 'Zapotlán De Juárez'
$ws.Range("B435").Value = 'Ahualulco De Mercado'
$ws.Range("B438").Value = 'Atotonilco El Alto'
$ws.Range("B439").Value = 'Autlán De Navarro'
$ws.Range("B445").Value = 'Encarnación De Díaz'
$ws.Range("B447").Value = 'Ixtlahuacán Del Río'
$ws.Range("B449").Value = 'Jilotlán De Los Dolores'
$ws.Range("B452").Value = 'Lagos De Moreno'
$ws.Range("B455").Value = 'Ojuelos De Jalisco'
$ws.Range("B463").Value = 'Talpa De Allende'
$ws.Range("B464").Value = 'Tamazula De Gordiano'
$ws.Range("B468").Value = 'Tepatitlán De Morelos'
$ws.Range("B470").Value = 'Tizapán El Alto'
$ws.Range("B471").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B476").Value = 'Unión De San Antonio'
$ws.Range("B477").Value = 'Valle De Juárez'
$ws.Range("B479").Value = 'Yahualica De González Gallo'
$ws.Range("B480").Value = 'Zacoalco De Torres'
$ws.Range("B482").Value = 'Zapotlán El Grande'
$ws.Range("B501").Value = 'Cojumatlán De Régules'
$ws.Range("B571").Value = 'Coatlán Del Río'
$ws.Range("B583").Value = 'Puente De Ixtla'
$ws.Range("B589").Value = 'Tetela Del Volcán'
$ws.Range("B591").Value = 'Tlaltizapán De Zapata'
$ws.Range("B599").Value = 'Zacualpan De Amilpas'
$ws.Range("B613").Value = 'Mier Y Noriega'
$ws.Range("B618").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B621").Value = 'Ayoquezco De Aldama'
$ws.Range("B625").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B626").Value = 'Ciénega De Zimatlán'
$ws.Range("B628").Value = 'Coicoyán De Las Flores'
$ws.Range("B629").Value = 'Constancia Del Rosario'
$ws.Range("B632").Value = 'Cuilápam De Guerrero'
$ws.Range("B633").Value = 'Eloxochitlán De Flores Magón'
$ws.Range("B634").Value = 'Fresnillo De Trujano'
$ws.Range("B635").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B636").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B637").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B639").Value = 'Huautla De Jiménez'
$ws.Range("B641").Value = 'Ixtlán De Juárez'
$ws.Range("B642").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B653").Value = 'Mariscala De Juárez'
$ws.Range("B654").Value = 'Mártires De Tacubaya'
$ws.Range("B656").Value = 'Mazatlán Villa De Flores'
$ws.Range("B658").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B659").Value = 'Mixistlán De La Reforma'
$ws.Range("B660").Value = 'Nejapa De Madero'
$ws.Range("B661").Value = 'Oaxaca De Juárez'
$ws.Range("B662").Value = 'Ocotlán De Morelos'
$ws.Range("B663").Value = 'Pinotepa De Don Luis'
$ws.Range("B664").Value = 'Putla Villa De Guerrero'
$ws.Range("B665").Value = 'Reforma De Pineda'
$ws.Range("B670").Value = 'San Agustín De Las Juntas'
$ws.Range("B683").Value = 'San Antonino El Alto'
$ws.Range("B686").Value = 'San Antonio De La Cal'
$ws.Range("B691").Value = 'San Baltazar Yatzachi El Bajo'
$ws.Range("B721").Value = 'San José Del Progreso'
$ws.Range("B729").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B739").Value = 'San Juan Del Estado'
$ws.Range("B740").Value = 'San Juan Del Río'
$ws.Range("B770").Value = 'San Martín De Los Cansecos'
$ws.Range("B784").Value = 'San Miguel Del Puerto'
$ws.Range("B785").Value = 'San Miguel El Grande'
$ws.Range("B802").Value = 'San Pablo Villa De Mitla'
$ws.Range("B808").Value = 'San Pedro El Alto'
$ws.Range("B821").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B822").Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range("B834").Value = 'Santa Ana Del Valle'
$ws.Range("B849").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B854").Value = 'Santa Inés De Zaragoza'
$ws.Range("B855").Value = 'Santa Inés Del Monte'
$ws.Range("B857").Value = 'Santa Lucía Del Camino'
$ws.Range("B866").Value = 'Santa María Del Tule'
$ws.Range("B873").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B929").Value = 'Santo Domingo De Morelos'
$ws.Range("B944").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B945").Value = 'Tataltepec De Valdés'
$ws.Range("B946").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B947").Value = 'Teotitlán De Flores Magón'
$ws.Range("B948").Value = 'Teotitlán Del Valle'
$ws.Range("B951").Value = 'Tlacolula De Matamoros'
$ws.Range("B953").Value = 'Tlalixtac De Cabrera'
$ws.Range("B954").Value = 'Totontepec Villa De Morelos'
$ws.Range("B957").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B958").Value = 'Villa De Etla'
$ws.Range("B959").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B960").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B961").Value = 'Villa De Zaachila'
$ws.Range("B963").Value = 'Villa Sola De Vega'
$ws.Range("B964").Value = 'Villa Talea De Castro'
$ws.Range("B965").Value = 'Villa Tejúpam De La Unión'
$ws.Range("B968").Value = 'Yutanduchi De Guerrero'
$ws.Range("B972").Value = 'Zimatlán De Álvarez'
$ws.Range("B998").Value = 'Ayotoxco De Guerrero'
$ws.Range("B1003").Value = 'Chalchicomula De Sesma'
$ws.Range("B1013").Value = 'Chila De La Sal'
$ws.Range("B1027").Value = 'Cuayuca De Andrade'
$ws.Range("B1028").Value = 'Cuetzalan Del Progreso'
$ws.Range("B1043").Value = 'Huehuetlán El Chico'
$ws.Range("B1044").Value = 'Huehuetlán El Grande'
$ws.Range("B1049").Value = 'Huitzilan De Serdán'
$ws.Range("B1051").Value = 'Ixcamilpa De Guerrero'
$ws.Range("B1055").Value = 'Izúcar De Matamoros'
$ws.Range("B1065").Value = 'Los Reyes De Juárez'
$ws.Range("B1066").Value = 'Mazapiltepec De Juárez'
$ws.Range("B1079").Value = 'Palmar De Bravo'
$ws.Range("B1089").Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Range("B1104").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B1108").Value = 'San Salvador El Seco'
$ws.Range("B1109").Value = 'San Salvador El Verde'
$ws.Range("B1118").Value = 'Tecali De Herrera'
$ws.Range("B1126").Value = 'Tepanco De López'
$ws.Range("B1127").Value = 'Tepango De Rodríguez'
$ws.Range("B1128").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B1134").Value = 'Tepexi De Rodríguez'
$ws.Range("B1136").Value = 'Tetela De Ocampo'
$ws.Range("B1137").Value = 'Teteles De Avila Castillo'
$ws.Range("B1142").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B1154").Value = 'Totoltepec De Guerrero'
$ws.Range("B1156").Value = 'Tuzamapan De Galeana'
$ws.Range("B1160").Value = 'Xayacatlán De Bravo'
$ws.Range("B1166").Value = 'Xochitlán De Vicente Suárez'
$ws.Range("B1182").Value = 'Amealco De Bonfil'
$ws.Range("B1184").Value = 'Cadereyta De Montes'
$ws.Range("B1186").Value = 'Jalpan De Serra'
$ws.Range("B1187").Value = 'Pinal De Amoles'
$ws.Range("B1190").Value = 'San Juan Del Río'
$ws.Range("B1197").Value = 'Axtla De Terrazas'
$ws.Range("B1202").Value = 'Ciudad Del Maíz'
$ws.Range("B1210").Value = 'Santa María Del Río'
$ws.Range("B1218").Value = 'Villa De La Paz'
$ws.Range("B1219").Value = 'Villa De Ramos'
$ws.Range("B1220").Value = 'Villa De Reyes'
$ws.Range("B1252").Value = 'Jalpa De Méndez'
$ws.Range("B1270").Value = 'Soto La Marina'
$ws.Range("B1276").Value = 'Acuamanala De Miguel Hidalgo'
$ws.Range("B1278").Value = 'Amaxac De Guerrero'
$ws.Range("B1279").Value = 'Apetatitlán De Antonio Carvajal'
$ws.Range("B1284").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1291").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1294").Value = 'Muñoz De Domingo Arenas'
$ws.Range("B1295").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1298").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B1304").Value = 'San Pablo Del Monte'
$ws.Range("B1312").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1315").Value = 'Tetla De La Solidaridad'
$ws.Range("B1327").Value = 'Ziltlaltépec De Trinidad Sánchez Santos'
$ws.Range("B1337").Value = 'Amatlán De Los Reyes'
$ws.Range("B1344").Value = 'Boca Del Río'
$ws.Range("B1349").Value = 'Castillo De Teayo'
$ws.Range("B1351").Value = 'Cazones De Herrera'
$ws.Range("B1364").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1365").Value = 'Cosautlán De Carvajal'
$ws.Range("B1379").Value = 'Hueyapan De Ocampo'
$ws.Range("B1380").Value = 'Ignacio De La Llave'
$ws.Range("B1384").Value = 'Ixhuacán De Los Reyes'
$ws.Range("B1385").Value = 'Ixhuatlán De Madero'
$ws.Range("B1386").Value = 'Ixhuatlán Del Café'
$ws.Range("B1394").Value = 'Juchique De Ferrer'
$ws.Range("B1399").Value = 'Las Vigas De Ramírez'
$ws.Range("B1400").Value = 'Lerdo De Tejada'
$ws.Range("B1404").Value = 'Martínez De La Torre'
$ws.Range("B1407").Value = 'Medellín De Bravo'
$ws.Range("B1411").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1420").Value = 'Paso De Ovejas'
$ws.Range("B1421").Value = 'Paso Del Macho'
$ws.Range("B1425").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1433").Value = 'Sayula De Alemán'
$ws.Range("B1436").Value = 'Soledad De Doblado'
$ws.Range("B1460").Value = 'Tlacotepec De Mejía'
$ws.Range("B1466").Value = 'Vega De Alatorre'
$ws.Range("B1478").Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range("B1479").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1490").Value = 'Jiménez Del Teul'
$ws.Range("B1495").Value = 'Nochistlán De Mejía'
$ws.Range("B1496").Value = 'Noria De Ángeles'

# --- Minor floating point precision updates (percentage recalculation) ---
$ws.Range("D209").Value = 0.0009609481354936872
$ws.Range("D354").Value = 0.0009609481354936872
$ws.Range("D1035").Value = 0.0009609481354936872
$ws.Range("D1079").Value = 0.0009609481354936872
$ws.Range("D1109").Value = 0.0009609481354936872
$ws.Range("D1282").Value = 0.0009609481354936872
$ws.Range("D1504").Value = 0.0009609481354936872
$ws.Range("D283").Value = 0.009048928275898888
$ws.Range("D1120").Value = 0.009048928275898888

# --- Remove trailing metadata/footer rows (1506:1511), shrinking used range to A1:D1505 ---
$ws.Rows("1506:1511").Delete()

